# Auto-generated script applying 210 cell value updates
# to match the target diff for Jogos_do_Dia_Betfair_Back_Lay_2025-11-21.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.26
$ws.Range("G2").Value = 2.32
$ws.Range("AO2").Value = 1000
$ws.Range("S3").Value = 2.84
$ws.Range("W3").Value = 1.48
$ws.Range("H5").Value = 1.18
$ws.Range("I6").Value = 9.199999999999999
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 2.54
$ws.Range("H7").Value = 1.91
$ws.Range("I7").Value = 2.56
$ws.Range("AF11").Value = 11
$ws.Range("F12").Value = 2.32
$ws.Range("G12").Value = 2.62
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 3.25
$ws.Range("P12").Value = 1.68
$ws.Range("Q12").Value = 1.99
$ws.Range("F13").Value = 2.04
$ws.Range("G13").Value = 2.28
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 3.25
$ws.Range("K13").Value = 3.8
$ws.Range("Q13").Value = 1.9
$ws.Range("G15").Value = 2.9
$ws.Range("I15").Value = 3.05
$ws.Range("W15").Value = 1.53
$ws.Range("AB15").Value = 1000
$ws.Range("AE15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("AN16").Value = 23
$ws.Range("Q17").Value = 1.92
$ws.Range("I18").Value = 6.4
$ws.Range("P18").Value = 2.32
$ws.Range("S18").Value = 2.7
$ws.Range("AF18").Value = 11
$ws.Range("AN20").Value = 19.5
$ws.Range("I21").Value = 5.5
$ws.Range("F22").Value = 3.5
$ws.Range("G22").Value = 4
$ws.Range("J22").Value = 3.7
$ws.Range("N22").Value = 5.3
$ws.Range("P22").Value = 2.42
$ws.Range("S22").Value = 2.5
$ws.Range("X22").Value = 28
$ws.Range("Y22").Value = 14.5
$ws.Range("Z22").Value = 17
$ws.Range("AA22").Value = 28
$ws.Range("AB22").Value = 24
$ws.Range("AC22").Value = 9.6
$ws.Range("AD22").Value = 12
$ws.Range("AE22").Value = 21
$ws.Range("AF22").Value = 1000
$ws.Range("AG22").Value = 16.5
$ws.Range("AH22").Value = 15.5
$ws.Range("AI22").Value = 30
$ws.Range("AK22").Value = 1000
$ws.Range("AN22").Value = 28
$ws.Range("AO22").Value = 11
$ws.Range("G24").Value = 7.8
$ws.Range("P24").Value = 2.58
$ws.Range("Q24").Value = 1.57
$ws.Range("R24").Value = 1.61
$ws.Range("S24").Value = 2.38
$ws.Range("Q25").Value = 2.36
$ws.Range("F26").Value = 1.8
$ws.Range("N26").Value = 2.6
$ws.Range("P26").Value = 1.76
$ws.Range("Q26").Value = 1.81
$ws.Range("S27").Value = 3.9
$ws.Range("T27").Value = 1.91
$ws.Range("U27").Value = 1.9
$ws.Range("F28").Value = 2.46
$ws.Range("Q29").Value = 1.8
$ws.Range("S29").Value = 3
$ws.Range("T29").Value = 1.67
$ws.Range("U29").Value = 2.22
$ws.Range("P30").Value = 1.55
$ws.Range("Q31").Value = 1.9
$ws.Range("P32").Value = 1.66
$ws.Range("Q32").Value = 1.93
$ws.Range("H33").Value = 3.35
$ws.Range("L33").Value = 1.37
$ws.Range("M33").Value = 1.07
$ws.Range("Q33").Value = 1.9
$ws.Range("S33").Value = 3.35
$ws.Range("AB33").Value = 10.5
$ws.Range("AC33").Value = 9.4
$ws.Range("AG33").Value = 13
$ws.Range("F35").Value = 1.44
$ws.Range("S35").Value = 2.14
$ws.Range("M37").Value = 1.07
$ws.Range("R37").Value = 1.29
$ws.Range("S37").Value = 3.85
$ws.Range("X37").Value = 15
$ws.Range("Z37").Value = 20
$ws.Range("AA37").Value = 48
$ws.Range("AC37").Value = 7.8
$ws.Range("AD37").Value = 15
$ws.Range("AE37").Value = 38
$ws.Range("AF37").Value = 22
$ws.Range("AG37").Value = 17
$ws.Range("AH37").Value = 23
$ws.Range("AI37").Value = 55
$ws.Range("AK37").Value = 48
$ws.Range("AN37").Value = 48
$ws.Range("AO37").Value = 34
$ws.Range("I38").Value = 3.65
$ws.Range("Y38").Value = 9
$ws.Range("F39").Value = 3.2
$ws.Range("M39").Value = 1.09
$ws.Range("R39").Value = 1.22
$ws.Range("S39").Value = 4.7
$ws.Range("U39").Value = 1.86
$ws.Range("M42").Value = 1.08
$ws.Range("O42").Value = 1.35
$ws.Range("P42").Value = 1.59
$ws.Range("R42").Value = 1.24
$ws.Range("U42").Value = 1.78
$ws.Range("V42").Value = 1.47
$ws.Range("I43").Value = 1.92
$ws.Range("K43").Value = 4.8
$ws.Range("P43").Value = 2.56
$ws.Range("Q43").Value = 1.51
$ws.Range("R43").Value = 1.63
$ws.Range("U43").Value = 2.44
$ws.Range("V43").Value = 1.91
$ws.Range("Y43").Value = 16.5
$ws.Range("AA43").Value = 25
$ws.Range("AD43").Value = 13
$ws.Range("Q44").Value = 1.34
$ws.Range("S44").Value = 1.94
$ws.Range("S45").Value = 2.32
$ws.Range("G46").Value = 1.99
$ws.Range("W46").Value = 2
$ws.Range("G47").Value = 2.28
$ws.Range("N47").Value = 4.9
$ws.Range("O47").Value = 1.21
$ws.Range("R47").Value = 1.53
$ws.Range("S47").Value = 2.52
$ws.Range("T47").Value = 1.58
$ws.Range("U47").Value = 2.44
$ws.Range("W47").Value = 1.78
$ws.Range("X47").Value = 27
$ws.Range("Y47").Value = 21
$ws.Range("Z47").Value = 34
$ws.Range("AB47").Value = 16.5
$ws.Range("AC47").Value = 11.5
$ws.Range("AD47").Value = 18
$ws.Range("AE47").Value = 42
$ws.Range("AF47").Value = 20
$ws.Range("AG47").Value = 14
$ws.Range("AH47").Value = 19
$ws.Range("AI47").Value = 46
$ws.Range("AJ47").Value = 34
$ws.Range("AK47").Value = 26
$ws.Range("AL47").Value = 36
$ws.Range("AN47").Value = 14.5
$ws.Range("AO47").Value = 30
$ws.Range("N48").Value = 1.1
$ws.Range("I49").Value = 2.18
$ws.Range("V49").Value = 1.84
$ws.Range("G50").Value = 2.84
$ws.Range("J50").Value = 3.3
$ws.Range("G51").Value = 5.8
$ws.Range("F52").Value = 1.73
$ws.Range("K52").Value = 3.95
$ws.Range("R52").Value = 1.32
$ws.Range("Q55").Value = 1.82
$ws.Range("W55").Value = 1.44
$ws.Range("F56").Value = 1.76
$ws.Range("G56").Value = 1.81
$ws.Range("T56").Value = 1.87
$ws.Range("W56").Value = 2.22
$ws.Range("AB56").Value = 7.8
$ws.Range("AC56").Value = 8.6
$ws.Range("AD56").Value = 1000
$ws.Range("AF56").Value = 10.5
$ws.Range("AG56").Value = 10.5
$ws.Range("AK56").Value = 1000
$ws.Range("AN56").Value = 15
$ws.Range("W58").Value = 1.83
$ws.Range("G59").Value = 3.85
$ws.Range("K60").Value = 4
$ws.Range("F61").Value = 3.75
$ws.Range("G61").Value = 5.9
$ws.Range("H61").Value = 1.77
$ws.Range("I61").Value = 2.24
$ws.Range("K61").Value = 1000
$ws.Range("P61").Value = 1.9
$ws.Range("V61").Value = 1.8
$ws.Range("W61").Value = 1.28
$ws.Range("G63").Value = 2.16
$ws.Range("H63").Value = 3.85
$ws.Range("I63").Value = 3.9
$ws.Range("W63").Value = 1.87
$ws.Range("Z63").Value = 26
$ws.Range("AB63").Value = 9.6
$ws.Range("AD63").Value = 15
$ws.Range("AE63").Value = 44
$ws.Range("Q64").Value = 1.51
$ws.Range("F65").Value = 1.53
$ws.Range("G65").Value = 1.71
$ws.Range("H65").Value = 4.1
$ws.Range("I65").Value = 7
$ws.Range("K65").Value = 500
$ws.Range("V65").Value = 1.16
$ws.Range("W65").Value = 2.4
